$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference a default/unstyled cell so that after forcing a cell to text format
# (needed so numeric-looking strings like "1.002" or "1.000" are preserved exactly,
# trailing zeros and all) we can restore the original "no explicit style" look,
# avoiding a stray style index on the written cells.
$normalStyle = $ws.Range("B2").Style

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $normalStyle
}

Set-TextValue $ws.Range("D2") "30.299.31"
Set-TextValue $ws.Range("E2") "  +0.12%  "

Set-TextValue $ws.Range("D3") "1.867.89"
Set-TextValue $ws.Range("E3") "  +0.10%  "

Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.14%  "

Set-TextValue $ws.Range("D5") "234.83"
Set-TextValue $ws.Range("E5") "  -1.11%  "

Set-TextValue $ws.Range("E6") "  +0.10%  "

Set-TextValue $ws.Range("D7") "0.4698"
Set-TextValue $ws.Range("E7") "  +0.24%  "

Set-TextValue $ws.Range("D8") "0.2867"
Set-TextValue $ws.Range("E8") "  -0.08%  "

Set-TextValue $ws.Range("D9") "0.06571"
Set-TextValue $ws.Range("E9") "  +0.37%  "

Set-TextValue $ws.Range("E10") "  -2.75%  "

Set-TextValue $ws.Range("D11") "0.08016"
Set-TextValue $ws.Range("E11") "  +1.56%  "

Set-TextValue $ws.Range("D12") "96.78"
Set-TextValue $ws.Range("E12") "  -1.09%  "

Set-TextValue $ws.Range("D13") "1.870.57"
Set-TextValue $ws.Range("E13") "  +0.24%  "

Set-TextValue $ws.Range("D14") "5.114"
Set-TextValue $ws.Range("E14") "  -1.30%  "

Set-TextValue $ws.Range("E15") "  +0.34%  "

Set-TextValue $ws.Range("D16") "269.04"
Set-TextValue $ws.Range("E16") "  -3.46%  "

Set-TextValue $ws.Range("D17") "30.323.66"
Set-TextValue $ws.Range("E17") "  +0.19%  "

Set-TextValue $ws.Range("D18") "14.00"
Set-TextValue $ws.Range("E18") "  +2.82%  "

Set-TextValue $ws.Range("D19") "0.000007602"
Set-TextValue $ws.Range("E19") "  +3.47%  "

Set-TextValue $ws.Range("E20") "  +0.02%  "

Set-TextValue $ws.Range("D21") "2.118.25"
Set-TextValue $ws.Range("E21") "  +0.25%  "

Set-TextValue $ws.Range("D22") "1.002"
Set-TextValue $ws.Range("E22") "  +0.17%  "

Set-TextValue $ws.Range("D23") "5.258"
Set-TextValue $ws.Range("E23") "  -2.53%  "

Set-TextValue $ws.Range("D24") "6.204"

Set-TextValue $ws.Range("D25") "9.384"
Set-TextValue $ws.Range("E25") "  +0.99%  "

Set-TextValue $ws.Range("D26") "167.98"
Set-TextValue $ws.Range("E26") "  -0.40%  "

Set-TextValue $ws.Range("E27") "  -1.14%  "

Set-TextValue $ws.Range("E28") "  +0.36%  "

Set-TextValue $ws.Range("E29") "  -1.03%  "

Set-TextValue $ws.Range("D30") "0.09870"
Set-TextValue $ws.Range("E30") "  +0.49%  "

Set-TextValue $ws.Range("D31") "4.358"
Set-TextValue $ws.Range("E31") "  -0.81%  "

Set-TextValue $ws.Range("D32") "1.462"
Set-TextValue $ws.Range("E32") "  -1.24%  "

Set-TextValue $ws.Range("D33") "4.060"

Set-TextValue $ws.Range("D34") "0.04713"
Set-TextValue $ws.Range("E34") "  -0.70%  "

Set-TextValue $ws.Range("D35") "1.135"
Set-TextValue $ws.Range("E35") "  -0.73%  "

Set-TextValue $ws.Range("D36") "0.6990"
Set-TextValue $ws.Range("E36") "  -1.39%  "

Set-TextValue $ws.Range("D37") "2.714"
Set-TextValue $ws.Range("E37") "  +0.22%  "

Set-TextValue $ws.Range("D38") "0.01873"
Set-TextValue $ws.Range("E38") "  -0.29%  "

Set-TextValue $ws.Range("D39") "2.673"
Set-TextValue $ws.Range("E39") "  +1.95%  "

Set-TextValue $ws.Range("D40") "6.267"
Set-TextValue $ws.Range("E40") "  -0.59%  "

Set-TextValue $ws.Range("D41") "71.78"
Set-TextValue $ws.Range("E41") "  -6.64%  "

Set-TextValue $ws.Range("D42") "1.955"
Set-TextValue $ws.Range("E42") "  -0.38%  "

Set-TextValue $ws.Range("D43") "0.8421"

# Rows 44 and 45 swap content: PaxDollar moves up to row 44,
# TheSandbox moves down to row 45 (with updated price/volume).
Set-TextValue $ws.Range("B44") "PaxDollar"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D44") "1.000"
Set-TextValue $ws.Range("E44") "  +0.04%  "

Set-TextValue $ws.Range("B45") "TheSandbox"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D45") "0.4159"
Set-TextValue $ws.Range("E45") "  -0.76%  "

Set-TextValue $ws.Range("D46") "102.76"
Set-TextValue $ws.Range("E46") "  -0.51%  "

Set-TextValue $ws.Range("D47") "7.046"
Set-TextValue $ws.Range("E47") "  -2.51%  "

Set-TextValue $ws.Range("D48") "9.137"
Set-TextValue $ws.Range("E48") "  -1.93%  "

Set-TextValue $ws.Range("D49") "910.22"
Set-TextValue $ws.Range("E49") "  -4.82%  "

Set-TextValue $ws.Range("E50") "  +0.39%  "

Set-TextValue $ws.Range("E51") "  +1.12%  "
